# NV-40 Sang sang 7-2024.xlsx
# "xoa cac phan tong cua ti le chiet khau"
#   - Set the discount-rate total (M3) on "Đơn sale chính" to 0
#   - Remove the "... tại HỆ THỐNG" rows (and their total row) from "Lương"

$wb = $excel.ActiveWorkbook

# 1) Sheet "Đơn sale chính": the Total row's discount-rate cell (M3) becomes 0
$wsOrders = $wb.Worksheets.Item("Đơn sale chính")
$wsOrders.Range("M3").Value = 0

# 2) Sheet "Lương": drop the 7 "... tại HỆ THỐNG" line items (rows 4-10)
$wsSalary = $wb.Worksheets.Item("Lương")
$wsSalary.Range("A4:B10").EntireRow.Delete()

# After the above deletion, "Tổng lương tại HỆ THỐNG" (originally row 35)
# has shifted up to row 28 - remove it too.
$wsSalary.Range("A28:B28").EntireRow.Delete()
